$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.2281368821292776
$ws.Range("C2").Value = 0.4942965779467681
$ws.Range("J2").Value = 0.01520912547528517
$ws.Range("P2").Value = 0.1444866920152091
$ws.Range("S2").Value = 0.1178707224334601

# Row 3
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.007575757575757576
$ws.Range("J3").Value = 0.04545454545454546
$ws.Range("P3").Value = 0.7348484848484849
$ws.Range("S3").Value = 0.196969696969697

# Row 4
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2916666666666667

# Row 6
$ws.Range("B6").Value = 0.05045871559633028
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("F6").Value = 0.07339449541284404
$ws.Range("J6").Value = 0.2155963302752294
$ws.Range("O6").Value = 0.01834862385321101
$ws.Range("Q6").Value = 0.2155963302752294
$ws.Range("R6").Value = 0.06422018348623854
$ws.Range("S6").Value = 0.3486238532110092

# Row 7
$ws.Range("B7").Value = 0.08
$ws.Range("D7").Value = 0.02
$ws.Range("F7").Value = 0.04666666666666667
$ws.Range("J7").Value = 0.1866666666666667
$ws.Range("O7").Value = 0.006666666666666667
$ws.Range("Q7").Value = 0.1266666666666667
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.4333333333333333

# Row 8
$ws.Range("B8").Value = 0.09043927648578812
$ws.Range("D8").Value = 0.0310077519379845
$ws.Range("F8").Value = 0.05167958656330749
$ws.Range("J8").Value = 0.08527131782945736
$ws.Range("O8").Value = 0.02325581395348837
$ws.Range("Q8").Value = 0.1757105943152455
$ws.Range("R8").Value = 0.07493540051679587
$ws.Range("S8").Value = 0.4677002583979328

# Row 9
$ws.Range("B9").Value = 0.06
$ws.Range("D9").Value = 0.025
$ws.Range("F9").Value = 0.065
$ws.Range("J9").Value = 0.13
$ws.Range("O9").Value = 0.005
$ws.Range("Q9").Value = 0.155
$ws.Range("R9").Value = 0.105
$ws.Range("S9").Value = 0.455

# Row 10
$ws.Range("B10").Value = 0.09991876523151909
$ws.Range("D10").Value = 0.02030869212022746
$ws.Range("E10").Value = 0.0008123476848090983
$ws.Range("F10").Value = 0.08204711616571893
$ws.Range("J10").Value = 0.1064175467099919
$ws.Range("O10").Value = 0.02030869212022746
$ws.Range("Q10").Value = 0.190089358245329
$ws.Range("R10").Value = 0.08285946385052803
$ws.Range("S10").Value = 0.3972380178716491

# Row 11
$ws.Range("G11").Value = 0.1265822784810127
$ws.Range("J11").Value = 0.08438818565400844
$ws.Range("K11").Value = 0.1772151898734177
$ws.Range("L11").Value = 0.6033755274261603
$ws.Range("S11").Value = 0.008438818565400843

# Row 12
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2312925170068027
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.02040816326530612
$ws.Range("S12").Value = 0.02040816326530612

# Row 13
$ws.Range("G13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.1176470588235294

# Row 15
$ws.Range("F15").Value = 0.02094240837696335
$ws.Range("H15").Value = 0.1361256544502618
$ws.Range("I15").Value = 0.07329842931937172
$ws.Range("J15").Value = 0.4083769633507853
$ws.Range("K15").Value = 0.05235602094240838
$ws.Range("M15").Value = 0.005235602094240838
$ws.Range("O15").Value = 0.0418848167539267
$ws.Range("S15").Value = 0.2617801047120419

# Row 16
$ws.Range("F16").Value = 0.01219512195121951
$ws.Range("H16").Value = 0.1463414634146341
$ws.Range("I16").Value = 0.09146341463414634
$ws.Range("J16").Value = 0.4085365853658536
$ws.Range("K16").Value = 0.06707317073170732
$ws.Range("M16").Value = 0.01829268292682927
$ws.Range("O16").Value = 0.0426829268292683
$ws.Range("S16").Value = 0.2134146341463415

# Row 17
$ws.Range("F17").Value = 0.02791878172588833
$ws.Range("H17").Value = 0.1218274111675127
$ws.Range("I17").Value = 0.1116751269035533
$ws.Range("J17").Value = 0.4720812182741117
$ws.Range("K17").Value = 0.05583756345177665
$ws.Range("M17").Value = 0.01269035532994924
$ws.Range("O17").Value = 0.06852791878172589
$ws.Range("S17").Value = 0.1294416243654822

# Row 18
$ws.Range("F18").Value = 0.0223463687150838
$ws.Range("H18").Value = 0.0893854748603352
$ws.Range("I18").Value = 0.111731843575419
$ws.Range("J18").Value = 0.5251396648044693
$ws.Range("K18").Value = 0.07262569832402235
$ws.Range("M18").Value = 0.00558659217877095
$ws.Range("O18").Value = 0.05027932960893855
$ws.Range("S18").Value = 0.1229050279329609

# Row 19
$ws.Range("F19").Value = 0.00996168582375479
$ws.Range("H19").Value = 0.210727969348659
$ws.Range("I19").Value = 0.08352490421455938
$ws.Range("J19").Value = 0.3701149425287356
$ws.Range("K19").Value = 0.1019157088122605
$ws.Range("M19").Value = 0.01839080459770115
$ws.Range("N19").Value = 0.0007662835249042146
$ws.Range("O19").Value = 0.06283524904214559
$ws.Range("S19").Value = 0.1417624521072797
